$wb = $excel.ActiveWorkbook

# --- Fig.2C (sheet 6) gets a reworked table layout ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Select()

# Header row
$ws6.Range("C1").Value = "Unit_Log"
$ws6.Range("D1").Value = "Annotation_label_1"
$ws6.Range("E1").Value = "Annotation_2"
$ws6.Range("F1").Value = "Annotation_label_2"
$ws6.Range("G1").Value = "Annotation_1"

# Row 2 keeps C/D but gets new E/F values, D is relabeled to LPS
$ws6.Range("D2").Value = "LPS"
$ws6.Range("E2").Value = "Time"
$ws6.Range("F2").Value = "0h"

# Rows 3-9: clear C/D, set F to 0h (rows3-5) / 6h (rows6-9)
$ws6.Range("C3").Value = ""
$ws6.Range("D3").Value = ""
$ws6.Range("F3").Value = "0h"

$ws6.Range("C4").Value = ""
$ws6.Range("D4").Value = ""
$ws6.Range("F4").Value = "0h"

$ws6.Range("C5").Value = ""
$ws6.Range("D5").Value = ""
$ws6.Range("F5").Value = "0h"

$ws6.Range("C6").Value = ""
$ws6.Range("D6").Value = ""
$ws6.Range("F6").Value = "6h"

$ws6.Range("C7").Value = ""
$ws6.Range("D7").Value = ""
$ws6.Range("F7").Value = "6h"

$ws6.Range("C8").Value = ""
$ws6.Range("D8").Value = ""
$ws6.Range("F8").Value = "6h"

$ws6.Range("C9").Value = ""
$ws6.Range("D9").Value = ""
$ws6.Range("F9").Value = "6h"

# Rows 10-13: clear C/D, set E to "+", F to 0h
$ws6.Range("C10").Value = ""
$ws6.Range("D10").Value = ""
$ws6.Range("E10").Value = "+"
$ws6.Range("F10").Value = "0h"

$ws6.Range("C11").Value = ""
$ws6.Range("D11").Value = ""
$ws6.Range("E11").Value = "+"
$ws6.Range("F11").Value = "0h"

$ws6.Range("C12").Value = ""
$ws6.Range("D12").Value = ""
$ws6.Range("E12").Value = "+"
$ws6.Range("F12").Value = "0h"

$ws6.Range("C13").Value = ""
$ws6.Range("D13").Value = ""
$ws6.Range("E13").Value = "'+"
$ws6.Range("F13").Value = "0h"

# Rows 14-17: clear C/D, set E to "+", F to 6h
$ws6.Range("C14").Value = ""
$ws6.Range("D14").Value = ""
$ws6.Range("E14").Value = "+"
$ws6.Range("F14").Value = "6h"

$ws6.Range("C15").Value = ""
$ws6.Range("D15").Value = ""
$ws6.Range("E15").Value = "+"
$ws6.Range("F15").Value = "6h"

$ws6.Range("C16").Value = ""
$ws6.Range("D16").Value = ""
$ws6.Range("E16").Value = "+"
$ws6.Range("F16").Value = "6h"

$ws6.Range("C17").Value = ""
$ws6.Range("D17").Value = ""
$ws6.Range("E17").Value = "+"
$ws6.Range("F17").Value = "6h"

$ws6.Range("B10").Select()

# --- restore the leftover selection state on Fig.1A ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Select()
$ws1.Range("A2:D17").Select()

# --- finally leave Fig.2C as the active tab/sheet ---
$ws6.Select()
$ws6.Range("B10").Select()
